$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New phrases to append ("Luận cách cục 12 cung" - Phụ Mẫu palace aspects)
$values = @(
    "Sát Phá Tham hội chiếu tại Phụ Mẫu",
    "Tử Phủ Vũ Tướng hội chiếu tại Phụ Mẫu",
    "Cơ Nguyệt Đồng Lương hội chiếu tại Phụ Mẫu",
    "Cự Nhật hội chiếu tại Phụ Mẫu",
    "Kình Đà hội chiếu tại Phụ Mẫu",
    "Xương Khúc hội chiếu tại Phụ Mẫu",
    "Hoả Linh hội chiếu tại Phụ Mẫu",
    "Không Kiếp hội chiếu tại Phụ Mẫu",
    "Quang Quý hội chiếu tại Phụ Mẫu",
    "Tả Hữu hội chiếu tại Phụ Mẫu",
    "Song Hao hội chiếu tại Phụ Mẫu",
    "Tang Hổ hội chiếu tại Phụ Mẫu",
    "Khốc Hư hội chiếu tại Phụ Mẫu",
    "Hình Riêu hội chiếu tại Phụ Mẫu",
    "Thai Toạ hội chiếu tại Phụ Mẫu",
    "Đào Hồng hội chiếu tại Phụ Mẫu",
    "Ấn Phù hội chiếu tại Phụ Mẫu",
    "Song Hao Quyền Lộc Kiếp Hoả hội chiếu tại Phụ Mẫu",
    "Tử Phủ Vũ Tướng Xương Khúc Khôi Việt Tả Hữu Khoa Quyền Lộc Long hội chiếu tại Phụ Mẫu",
    "Tử Phủ Vũ Tướng Tả Hữu Khoa Quyền Lộc Long Phượng hội chiếu tại Phụ Mẫu",
    "Tử Khúc Phá Dương Đà hội chiếu tại Phụ Mẫu",
    "Cơ Nguyệt Đồng Lương gặp Xương Khúc Tả Hữu hội chiếu tại Phụ Mẫu",
    "Cơ Nguyệt Đồng Lương Khoa Tả Hữu Quang Quý Quan Phúc hội chiếu tại Phụ Mẫu",
    "Sát Quyền hội chiếu tại Phụ Mẫu",
    "Lộc Mã hội chiếu tại Phụ Mẫu",
    "Kiếp Hư Hao Quyền hội chiếu tại Phụ Mẫu",
    "Tuế Hổ Phù Xương Khúc hội chiếu tại Phụ Mẫu",
    "Xương Khúc Tấu Long Phượng hội chiếu tại Phụ Mẫu",
    "Đào Hồng Riêu Tấu Cơ Vũ hội chiếu tại Phụ Mẫu",
    "Binh Hình Tướng Ấn hội chiếu tại Phụ Mẫu",
    "Hổ Tấu hội chiếu tại Phụ Mẫu",
    "Hình Riêu Y hội chiếu tại Phụ Mẫu",
    "Mã Hỏa Linh hội chiếu tại Phụ Mẫu",
    "Thai Tọa Hồng Đào hội chiếu tại Phụ Mẫu",
    "Tả Hữu Không Kiếp hội chiếu tại Phụ Mẫu",
    "Tả Hữu Binh Tướng hội chiếu tại Phụ Mẫu",
    "Đào Quyền hội chiếu tại Phụ Mẫu",
    "Đào Hồng Tả Cái Hữu hội chiếu tại Phụ Mẫu",
    "Quan Phúc Quang Tấu hội chiếu tại Phụ Mẫu",
    "Đào Hồng Xương Khúc hội chiếu tại Phụ Mẫu",
    "Đào Tử Phủ hội chiếu tại Phụ Mẫu",
    "Hổ Kình Sát hội chiếu tại Phụ Mẫu",
    "Hổ Tang Kiếp hội chiếu tại Phụ Mẫu",
    "Hổ Tang Không Kiếp hội chiếu tại Phụ Mẫu",
    "Âm Dương Lương hội chiếu tại Phụ Mẫu",
    "Cơ Lương Gia Hội hội chiếu tại Phụ Mẫu",
    "Nhật Chiếu Lôi Môn hội chiếu tại Phụ Mẫu",
    "Tả Hữu Xương Khúc hội chiếu tại Phụ Mẫu",
    "Tham Linh Triều Viên hội chiếu tại Phụ Mẫu",
    "Cự Hỏa Linh hội chiếu tại Phụ Mẫu",
    "Hồng Đào Kỵ hội chiếu tại Phụ Mẫu",
    "Cự Đồng Hình hội chiếu tại Phụ Mẫu",
    "Lương Phá hội chiếu tại Phụ Mẫu",
    "Khôi Việt hội chiếu tại Phụ Mẫu",
    "Tham Vũ Hỏa hội chiếu tại Phụ Mẫu",
    "Cự Kỵ hội chiếu tại Phụ Mẫu",
    "Kình Đà Hỏa hội chiếu tại Phụ Mẫu",
    "Kình Đà Không Kiếp hội chiếu tại Phụ Mẫu",
    "Đào Hồng Hỷ hội chiếu tại Phụ Mẫu",
    "Đào Hồng Kiếp Sát hội chiếu tại Phụ Mẫu",
    "Đào Hồng Kỵ hội chiếu tại Phụ Mẫu",
    "Đào Hồng Tả Phù Hữu Bật hội chiếu tại Phụ Mẫu",
    "Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc hội chiếu tại Phụ Mẫu",
    "Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc Long Trì Phượng Các hội chiếu tại Phụ Mẫu",
    "Lưu Hà Kiếp Sát hội chiếu tại Phụ Mẫu",
    "Phục Không Kiếp hội chiếu tại Phụ Mẫu",
    "Xương Khúc Khôi Việt hội chiếu tại Phụ Mẫu",
    "Khoa Quyền hội chiếu tại Phụ Mẫu",
    "Lộc Quyền hội chiếu tại Phụ Mẫu",
    "Tướng Binh Đào Hồng hội chiếu tại Phụ Mẫu",
    "Lương Khốc Tuế hội chiếu tại Phụ Mẫu",
    "Tướng Binh hội chiếu tại Phụ Mẫu",
    "Cự Tang hội chiếu tại Phụ Mẫu",
    "Cự Hỏa hội chiếu tại Phụ Mẫu",
    "Phá Hình Kỵ hội chiếu tại Phụ Mẫu",
    "Tang Trực Tuế hội chiếu tại Phụ Mẫu",
    "Tuế Xương Khúc hội chiếu tại Phụ Mẫu",
    "Tử Vi Tả Hữu hội chiếu tại Phụ Mẫu",
    "Tử Vi Tang Tả Hữu hội chiếu tại Phụ Mẫu"
)

$n = $values.Count
$startRow = 4251
$endRow = $startRow + $n - 1

$rng = $ws.Range("A" + $startRow + ":B" + $endRow)
$arr = New-Object 'object[,]' $n,2
for ($i = 0; $i -lt $n; $i++) {
    $arr[$i,0] = $values[$i]
    $arr[$i,1] = $values[$i]
}
$rng.Value = $arr

# Highlight-duplicates conditional format over the whole of column A,
# matching the standard "Highlight Duplicate Values" style already used
# elsewhere in the workbook (dark red text on light red fill).
$dupRange = $ws.Range("A1:A1048576")
$fc = $dupRange.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615
$fc.SetFirstPriority()

# Scroll/selection bookkeeping to match where the author ended up editing.
$ws.Range("D4336").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4310
$win.ScrollColumn = 1

Write-Output "ok"
